# Commit: "added conjtrols to select file. read. clear."
# The DigitalOutput06 row (row 7) is removed/cleared: Location (A7), Type (B7)
# and Card (D7) are cleared, leaving only the two already-blank BR Name / Eplan
# name cells (C7, E7). This also drops the now-unused "DigitalOutput06" shared
# string. Finally move the active selection from I29 to H18 (and let the view
# scroll back up instead of being pinned at row 19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").ClearContents() | Out-Null
$ws.Range("B7").ClearContents() | Out-Null
$ws.Range("D7").ClearContents() | Out-Null

$ws.Range("H18").Select() | Out-Null
